# Apply the edits described by the diff:
#  1. Results sheet (sheet 1), row 11: update Awarded Supplier info to reflect
#     a transition from Facility 10 moving to supplier C instead of B.
#  2. LP Model sheet (sheet 3), cell A2: the LP-format text gets a new
#     "Rule_0" constraint (limiting transitions to 4) and T_1_B / T_1_C are
#     added to the Binaries section.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update "Results" sheet, row 11
# ---------------------------------------------------------------------
$results = $wb.Worksheets.Item("Results")

$results.Range("G11").Value = "C"
$results.Range("H11").Value = 15
# Prefix with an apostrophe so Excel stores these as literal text (matching
# the original "3%"/"5%" text cells) instead of auto-converting them into
# percentage-formatted numbers, then clear the resulting "quote prefix"
# formatting so the cell style stays identical to the untouched cells.
$results.Range("I11").Value2 = "'4%"
$results.Range("I11").ClearFormats()
$results.Range("J11").Value = 14.4
$results.Range("K11").Value = 187.2
$results.Range("M11").Value = 2815.8
$results.Range("N11").Value2 = "'7%"
$results.Range("N11").ClearFormats()
$results.Range("O11").Value = 13.104

# ---------------------------------------------------------------------
# 2. Update "LP Model" sheet, cell A2 (the LP-format text blob)
# ---------------------------------------------------------------------
$lpModel = $wb.Worksheets.Item("LP Model")
$cell = $lpModel.Range("A2")
$text = $cell.Value2

# Insert the new Rule_0 constraint right before the Transition_10_A line,
# directly after the RebateTierUpper_C_1 constraint block.
$oldBlock1 = "RebateTierUpper_C_1: - 0.07 S_C + rebate_C + 97000000000 y_rebate_C_1`n <= 97000000000`nTransition_10_A:"
$newBlock1 = "RebateTierUpper_C_1: - 0.07 S_C + rebate_C + 97000000000 y_rebate_C_1`n <= 97000000000`nRule_0: T_10_A + T_10_B + T_1_B + T_1_C + T_2_A + T_2_C + T_3_A + T_3_B`n + T_4_A + T_4_B + T_5_A + T_5_B + T_6_A + T_6_B + T_7_A + T_7_B + T_8_A`n + T_8_B + T_9_A + T_9_B = 4`nTransition_10_A:"

$text = $text.Replace($oldBlock1, $newBlock1)

# Insert T_1_B and T_1_C into the Binaries list, right after T_10_B.
$oldBlock2 = "Binaries`nT_10_A`nT_10_B`nT_2_A"
$newBlock2 = "Binaries`nT_10_A`nT_10_B`nT_1_B`nT_1_C`nT_2_A"

$text = $text.Replace($oldBlock2, $newBlock2)

$cell.Value = $text
